$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 10-19 with the new scheme ordering / data (Gaussian-Quadrature moved up,
# three new Spiral-* rows inserted, followed by the previously-existing schemes).
# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value2 = 8
$ws.Range("B10").Value2 = 'Gaussian-Quadrature'
$ws.Range("C10").Value2 = 0.9997970327128718
$ws.Range("D10").Value2 = 0.9697660149721196
$ws.Range("E10").Value2 = 0.9997385364559017
$ws.Range("F10").Value2 = 0.9997970327128718
$ws.Range("G10").Value2 = 0.9683729531988541
$ws.Range("H10").Value2 = 1.001031999558077
$ws.Range("I10").Value2 = 0.9952941176470588
$ws.Range("J10").Value2 = 0.9697660149721196
$ws.Range("K10").Value2 = 0.9847522757140106
$ws.Range("L10").Value2 = 0.9922746542134413
$ws.Range("M10").Value2 = 0.9890001090908139

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value2 = 9
$ws.Range("B11").Value2 = 'Spiral-90deg-10rot-5space'
$ws.Range("C11").Value2 = 0.9927734059709842
$ws.Range("D11").Value2 = 0.9798195914221737
$ws.Range("E11").Value2 = 0.9947802115617808
$ws.Range("F11").Value2 = 0.9927734059709842
$ws.Range("G11").Value2 = 0.9835745737898508
$ws.Range("H11").Value2 = 1.000151500738855
$ws.Range("I11").Value2 = 0.9941399773397617
$ws.Range("J11").Value2 = 0.9798195914221737
$ws.Range("K11").Value2 = 0.9872999014919772
$ws.Range("L11").Value2 = 0.9900366537314808
$ws.Range("M11").Value2 = 0.9908732101372344

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = 'Spiral-90deg-15rot-5space'
$ws.Range("C12").Value2 = 0.992738379315072
$ws.Range("D12").Value2 = 0.9798645986460404
$ws.Range("E12").Value2 = 0.9947869749823363
$ws.Range("F12").Value2 = 0.992738379315072
$ws.Range("G12").Value2 = 0.9836817641425833
$ws.Range("H12").Value2 = 1.000159340623365
$ws.Range("I12").Value2 = 0.9941451741110581
$ws.Range("J12").Value2 = 0.9798645986460404
$ws.Range("K12").Value2 = 0.9873257868141883
$ws.Range("L12").Value2 = 0.9900320830646301
$ws.Range("M12").Value2 = 0.9908960386367424

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = 'Spiral-90deg-10rot-3space'
$ws.Range("C13").Value2 = 0.9927766999458111
$ws.Range("D13").Value2 = 0.9797980644512084
$ws.Range("E13").Value2 = 0.9948272211591419
$ws.Range("F13").Value2 = 0.9927766999458111
$ws.Range("G13").Value2 = 0.9836404378631705
$ws.Range("H13").Value2 = 1.000160864735405
$ws.Range("I13").Value2 = 0.994161525821143
$ws.Range("J13").Value2 = 0.9797980644512084
$ws.Range("K13").Value2 = 0.9873126428051752
$ws.Range("L13").Value2 = 0.9900446713754931
$ws.Range("M13").Value2 = 0.9908941356626467

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value2 = 12
$ws.Range("B14").Value2 = 'NoRotation-tilt60deg'
$ws.Range("C14").Value2 = 0.9968679999999996
$ws.Range("D14").Value2 = 0.9314519999999998
$ws.Range("E14").Value2 = 0.9981239999999998
$ws.Range("F14").Value2 = 0.9968679999999996
$ws.Range("G14").Value2 = 0.9341439999999996
$ws.Range("H14").Value2 = 1.053979999999997
$ws.Range("I14").Value2 = 0.9972920000000007
$ws.Range("J14").Value2 = 0.9314519999999998
$ws.Range("K14").Value2 = 0.9647879999999998
$ws.Range("L14").Value2 = 0.9808279999999997
$ws.Range("M14").Value2 = 0.9853099999999996

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value2 = 13
$ws.Range("B15").Value2 = 'Rotation-NoTilt'
$ws.Range("C15").Value2 = 1
$ws.Range("D15").Value2 = 0.89
$ws.Range("E15").Value2 = 1
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 0.89
$ws.Range("H15").Value2 = 1.1
$ws.Range("I15").Value2 = 1
$ws.Range("J15").Value2 = 0.89
$ws.Range("K15").Value2 = 0.9450000000000001
$ws.Range("L15").Value2 = 0.9725
$ws.Range("M15").Value2 = 0.9800000000000001

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = 'Rotation-60detTilt'
$ws.Range("C16").Value2 = 0.9968750080000031
$ws.Range("D16").Value2 = 0.9314582271999953
$ws.Range("E16").Value2 = 0.9981250047999984
$ws.Range("F16").Value2 = 0.9968750080000031
$ws.Range("G16").Value2 = 0.9341665536000002
$ws.Range("H16").Value2 = 1.053958451199999
$ws.Range("I16").Value2 = 0.9972916736000018
$ws.Range("J16").Value2 = 0.9314582271999953
$ws.Range("K16").Value2 = 0.9647916159999969
$ws.Range("L16").Value2 = 0.980833312
$ws.Range("M16").Value2 = 0.9853124863999997

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = 'HexGrid-90degTilt5degRes'
$ws.Range("C17").Value2 = 0.9905057259552789
$ws.Range("D17").Value2 = 0.992041522875699
$ws.Range("E17").Value2 = 0.9911008089428689
$ws.Range("F17").Value2 = 0.9905057259552789
$ws.Range("G17").Value2 = 0.9904970480807437
$ws.Range("H17").Value2 = 0.9900085496758524
$ws.Range("I17").Value2 = 0.9915071074709519
$ws.Range("J17").Value2 = 0.992041522875699
$ws.Range("K17").Value2 = 0.991571165909284
$ws.Range("L17").Value2 = 0.9910384459322814
$ws.Range("M17").Value2 = 0.9909434605002324

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = 'HexGrid-90degTilt22p5degRes'
$ws.Range("C18").Value2 = 0.9905469845476249
$ws.Range("D18").Value2 = 0.9982744380809565
$ws.Range("E18").Value2 = 0.9910885936214049
$ws.Range("F18").Value2 = 0.9905469845476249
$ws.Range("G18").Value2 = 0.9923287446189193
$ws.Range("H18").Value2 = 0.9846493877629974
$ws.Range("I18").Value2 = 0.9902983081668018
$ws.Range("J18").Value2 = 0.9982744380809565
$ws.Range("K18").Value2 = 0.9946815158511807
$ws.Range("L18").Value2 = 0.9926142501994029
$ws.Range("M18").Value2 = 0.991197742799784

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = 'HexGrid-60degTilt5degRes'
$ws.Range("C19").Value2 = 0.9888161938134855
$ws.Range("D19").Value2 = 1.00479870944407
$ws.Range("E19").Value2 = 0.9886897954349229
$ws.Range("F19").Value2 = 0.9888161938134855
$ws.Range("G19").Value2 = 0.9991932448620844
$ws.Range("H19").Value2 = 0.9801588706975204
$ws.Range("I19").Value2 = 0.9886335434724505
$ws.Range("J19").Value2 = 1.00479870944407
$ws.Range("K19").Value2 = 0.9967442524394963
$ws.Range("L19").Value2 = 0.9927802231264908
$ws.Range("M19").Value2 = 0.9917150596207555

# Rows 17-19 are brand new rows; give column A the same bold/centered/bordered style
# used by every other value in column A (rows 3-16), by copying the formatting from A16.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "done"